$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the format from the existing "sum" header (G1) into the new header
# cell so H1 picks up the same bold/bordered/centered style (s="1"),
# then set its text to "Save".
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New "Save" data column values (plain, unstyled numbers like the rest of
# the numeric columns).
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 0
